# Insert a new data row at row 292 (pushing the existing row 292..374 down
# to 293..375) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(292).Insert()

$ws.Range("A292").Value = 3
$ws.Range("B292").Value = 'Femacal de La Calera'
$ws.Range("C292").Value = 'Coquimbo'
$ws.Range("D292").Value = 44736
$ws.Range("E292").Value = 5
$ws.Range("F292").Value = 100112040
$ws.Range("G292").Value = 'Cilantro'
$ws.Range("H292").Value = 'Sin especificar'
$ws.Range("I292").Value = 'Primera'
$ws.Range("J292").Value = 160
$ws.Range("K292").Value = 3500
$ws.Range("L292").Value = 3500
$ws.Range("M292").Value = 3500
$ws.Range("N292").Value = '$/docena de atados (3 kilos)'
$ws.Range("O292").Value = 'Provincia de Quillota'
$ws.Range("P292").Value = 1167
$ws.Range("Q292").Value = 3
$ws.Range("R292").Value = 'Hortaliza'
